$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para5 = $tr.Paragraphs(5)
$para5.InsertAfter("`rTalk to me if interested`r")

# Now delete the paragraph mark at the end of para 6 (i.e. character at the very end of "Talk to me if interested")
$tr2 = $sh.TextFrame.TextRange
$p6 = $tr2.Paragraphs(6)
Write-Host "p6 length:" $p6.Length
# the paragraph mark is the character right after position Length in zero-index, try deleting Characters(Length+1,1) to remove pilcrow of p6 merging p6+p7
$c = $p6.Characters($p6.Length+1, 1)
Write-Host "c text: [" $c.Text "] len=" $c.Length
$c.Delete()

$tr3 = $sh.TextFrame.TextRange
Write-Host "count after merge:" $tr3.Paragraphs().Count
for ($i=1; $i -le $tr3.Paragraphs().Count; $i++) {
    Write-Host $i ": [" $tr3.Paragraphs($i).Text "]"
}
